$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.92745816989022
$ws.Range("C2").Value = 8.878632296859205
$ws.Range("D2").Value = 11.67128436937264
$ws.Range("F2").Value = 26.99214720608759
$ws.Range("G2").Value = 23.30384877354595
$ws.Range("H2").Value = 12.97674337305346
$ws.Range("I2").Value = 17.98293395393583
$ws.Range("J2").Value = 11.18012835938675
$ws.Range("M2").Value = 16.99044094625025
$ws.Range("O2").Value = 19.01279945314569
$ws.Range("B3").Value = 12.24468236580213
$ws.Range("C3").Value = 8.38114597930022
$ws.Range("D3").Value = 11.65138050341772
$ws.Range("F3").Value = 27.08882657028675
$ws.Range("G3").Value = 23.41569490584683
$ws.Range("H3").Value = 13.03556201886477
$ws.Range("I3").Value = 18.12221199475418
$ws.Range("J3").Value = 11.21838131080659
$ws.Range("M3").Value = 16.75020216162569
$ws.Range("O3").Value = 19.11276744880146
$ws.Range("B4").Value = 11.80471645433323
$ws.Range("C4").Value = 8.058956799081491
$ws.Range("D4").Value = 11.64106690495509
$ws.Range("F4").Value = 27.1563036661109
$ws.Range("G4").Value = 23.49511042110593
$ws.Range("H4").Value = 13.07419759309531
$ws.Range("I4").Value = 18.21223500421915
$ws.Range("J4").Value = 11.24364619219943
$ws.Range("M4").Value = 16.60251132205486
$ws.Range("O4").Value = 19.17938574375971
$ws.Range("B5").Value = 11.62036120828704
$ws.Range("C5").Value = 7.923524502004343
$ws.Range("D5").Value = 11.63734694334342
$ws.Range("F5").Value = 27.18583383998264
$ws.Range("G5").Value = 23.5301530517529
$ws.Range("H5").Value = 13.0905755579926
$ws.Range("I5").Value = 18.25005511412991
$ws.Range("J5").Value = 11.2543887888328
$ws.Range("M5").Value = 16.54234232437386
$ws.Range("O5").Value = 19.20784635036433
$ws.Range("B6").Value = 11.58944815478424
$ws.Range("C6").Value = 7.90078851673079
$ws.Range("D6").Value = 11.63675850233528
$ws.Range("F6").Value = 27.19085983097158
$ws.Range("G6").Value = 23.53613310003238
$ws.Range("H6").Value = 13.09333337252878
$ws.Range("I6").Value = 18.25640372549692
$ws.Range("J6").Value = 11.25619958396736
$ws.Range("M6").Value = 16.53235411617023
$ws.Range("O6").Value = 19.21265140671904
$ws.Range("B7").Value = 11.80225046659981
$ws.Range("C7").Value = 8.057146961353979
$ws.Range("D7").Value = 11.64101477687984
$ws.Range("F7").Value = 27.15669370091639
$ws.Range("G7").Value = 23.49557219472162
$ws.Range("H7").Value = 13.07441590665956
$ws.Range("I7").Value = 18.21274046118909
$ws.Range("J7").Value = 11.24378926090933
$ws.Range("M7").Value = 16.60169971364236
$ws.Range("O7").Value = 19.17976426189849
$ws.Range("B8").Value = 12.69642238663898
$ws.Range("C8").Value = 8.710608417005828
$ws.Range("D8").Value = 11.66402768043393
$ws.Range("F8").Value = 27.02379316844298
$ws.Range("G8").Value = 23.34017165570012
$ws.Range("H8").Value = 12.99650075561448
$ws.Range("I8").Value = 18.03002306053029
$ws.Range("J8").Value = 11.19294907924986
$ws.Range("M8").Value = 16.90768186568325
$ws.Range("O8").Value = 19.04617901804529
$ws.Range("B9").Value = 14.28011124703433
$ws.Range("C9").Value = 9.856964676308854
$ws.Range("D9").Value = 11.72413399656707
$ws.Range("F9").Value = 26.82790898720794
$ws.Range("G9").Value = 23.12154462358835
$ws.Range("H9").Value = 12.86371849996728
$ws.Range("I9").Value = 17.7073737730437
$ws.Range("J9").Value = 11.10735424716408
$ws.Range("M9").Value = 17.50355478417486
$ws.Range("O9").Value = 18.82595064448885
$ws.Range("B10").Value = 15.33468952391047
$ws.Range("C10").Value = 10.61453664811084
$ws.Range("D10").Value = 11.77719616804269
$ws.Range("F10").Value = 26.72389220825686
$ws.Range("G10").Value = 23.01454454667739
$ws.Range("H10").Value = 12.77837324204512
$ws.Range("I10").Value = 17.49193257207187
$ws.Range("J10").Value = 11.05306443624628
$ws.Range("M10").Value = 17.93526117822119
$ws.Range("O10").Value = 18.68982740839632
$ws.Range("B11").Value = 15.78998975504276
$ws.Range("C11").Value = 10.94048807507617
$ws.Range("D11").Value = 11.80321309004099
$ws.Range("F11").Value = 26.68531887119588
$ws.Range("G11").Value = 22.97772591275043
$ws.Range("H11").Value = 12.74220172907848
$ws.Range("I11").Value = 17.39858837364477
$ws.Range("J11").Value = 11.03023330939697
$ws.Range("M11").Value = 18.12961148119883
$ws.Range("O11").Value = 18.63352817414659
$ws.Range("B12").Value = 15.95883800885041
$ws.Range("C12").Value = 11.06121529429998
$ws.Range("D12").Value = 11.81332994886474
$ws.Range("F12").Value = 26.6719754881345
$ws.Range("G12").Value = 22.96550393625501
$ws.Range("H12").Value = 12.72888630255498
$ws.Range("I12").Value = 17.3639098047311
$ws.Range("J12").Value = 11.02185607718751
$ws.Range("M12").Value = 18.20285620827853
$ws.Range("O12").Value = 18.61302208289314
$ws.Range("B13").Value = 15.92263271850637
$ws.Range("C13").Value = 11.03533500521122
$ws.Range("D13").Value = 11.81113941275869
$ws.Range("F13").Value = 26.67479293980368
$ws.Range("G13").Value = 22.96805941047626
$ws.Range("H13").Value = 12.73173702250732
$ws.Range("I13").Value = 17.37134872920502
$ws.Range("J13").Value = 11.02364832351458
$ws.Range("M13").Value = 18.18709819333112
$ws.Range("O13").Value = 18.6174021976256
$ws.Range("B14").Value = 15.80395268245024
$ws.Range("C14").Value = 10.95047464550687
$ws.Range("D14").Value = 11.80404013445482
$ws.Range("F14").Value = 26.68419574801921
$ws.Range("G14").Value = 22.97668583998792
$ws.Range("H14").Value = 12.74109860322235
$ws.Range("I14").Value = 17.39572194988716
$ws.Range("J14").Value = 11.02953872903864
$ws.Range("M14").Value = 18.13564474321645
$ws.Range("O14").Value = 18.63182480213611
$ws.Range("B15").Value = 15.73079227411165
$ws.Range("C15").Value = 10.89814269028803
$ws.Range("D15").Value = 11.79972594471167
$ws.Range("F15").Value = 26.69011995438295
$ws.Range("G15").Value = 22.98219426270539
$ws.Range("H15").Value = 12.74688259739747
$ws.Range("I15").Value = 17.41073831527423
$ws.Range("J15").Value = 11.03318173329737
$ws.Range("M15").Value = 18.10408053206849
$ws.Range("O15").Value = 18.64076510040203
$ws.Range("B16").Value = 15.30443834283732
$ws.Range("C16").Value = 10.59285752077087
$ws.Range("D16").Value = 11.77553325557086
$ws.Range("F16").Value = 26.72658955910694
$ws.Range("G16").Value = 23.0171906727426
$ws.Range("H16").Value = 12.78079055378608
$ws.Range("I16").Value = 17.49812650615499
$ws.Range("J16").Value = 11.05459405720388
$ws.Range("M16").Value = 17.92251389554765
$ws.Range("O16").Value = 18.6936202480067
$ws.Range("B17").Value = 15.03658647499306
$ws.Range("C17").Value = 10.40077767823226
$ws.Range("D17").Value = 11.76116929881259
$ws.Range("F17").Value = 26.75120664481426
$ws.Range("G17").Value = 23.04170805624926
$ws.Range("H17").Value = 12.80227171009441
$ws.Range("I17").Value = 17.55292912397366
$ws.Range("J17").Value = 11.0682077508705
$ws.Range("M17").Value = 17.81056483983524
$ws.Range("O17").Value = 18.72748879278874
$ws.Range("B18").Value = 14.88022767122817
$ws.Range("C18").Value = 10.28854166501681
$ws.Range("D18").Value = 11.75308469957798
$ws.Range("F18").Value = 26.76618848242287
$ws.Range("G18").Value = 23.05692494109301
$ws.Range("H18").Value = 12.81487673391506
$ws.Range("I18").Value = 17.58488894313254
$ws.Range("J18").Value = 11.07621357527466
$ws.Range("M18").Value = 17.74598654120186
$ws.Range("O18").Value = 18.74749816520388
$ws.Range("B19").Value = 14.82689440234465
$ws.Range("C19").Value = 10.25023940854327
$ws.Range("D19").Value = 11.75037797966451
$ws.Range("F19").Value = 26.77140218696916
$ws.Range("G19").Value = 23.06226814302904
$ws.Range("H19").Value = 12.8191874373383
$ws.Range("I19").Value = 17.59578539757313
$ws.Range("J19").Value = 11.07895436269936
$ws.Range("M19").Value = 17.72409095085926
$ws.Range("O19").Value = 18.75436372777887
$ws.Range("B20").Value = 15.06533794829692
$ws.Range("C20").Value = 10.42140687570471
$ws.Range("D20").Value = 11.762680067914
$ws.Range("F20").Value = 26.74850091750965
$ws.Range("G20").Value = 23.03898260911762
$ws.Range("H20").Value = 12.79995916478779
$ws.Range("I20").Value = 17.54704988769277
$ws.Range("J20").Value = 11.06674037616465
$ws.Range("M20").Value = 17.82250191540295
$ws.Range("O20").Value = 18.72382863780362
$ws.Range("B21").Value = 15.83890896234934
$ws.Range("C21").Value = 10.97547368375679
$ws.Range("D21").Value = 11.80611822160124
$ws.Range("F21").Value = 26.68139957894362
$ws.Range("G21").Value = 22.97410524051635
$ws.Range("H21").Value = 12.73833851004383
$ws.Range("I21").Value = 17.38854480251345
$ws.Range("J21").Value = 11.02780128807609
$ws.Range("M21").Value = 18.15076786700119
$ws.Range("O21").Value = 18.6275664276721
$ws.Range("B22").Value = 16.32369219793208
$ws.Range("C22").Value = 11.32182523065465
$ws.Range("D22").Value = 11.83604842532909
$ws.Range("F22").Value = 26.64491200686575
$ws.Range("G22").Value = 22.94173821869909
$ws.Range("H22").Value = 12.70029234744181
$ws.Range("I22").Value = 17.28885055449681
$ws.Range("J22").Value = 11.00391697389457
$ws.Range("M22").Value = 18.36323403177831
$ws.Range("O22").Value = 18.56939556807368
$ws.Range("B23").Value = 16.06687051629937
$ws.Range("C23").Value = 11.13841799576927
$ws.Range("D23").Value = 11.8199349639254
$ws.Range("F23").Value = 26.66371018384209
$ws.Range("G23").Value = 22.95809027977729
$ws.Range("H23").Value = 12.72039439876917
$ws.Range("I23").Value = 17.3417030301787
$ws.Range("J23").Value = 11.0165212667478
$ws.Range("M23").Value = 18.25004553873293
$ws.Range("O23").Value = 18.60000707316005
$ws.Range("B24").Value = 15.05234677656864
$ws.Range("C24").Value = 10.41208604174223
$ws.Range("D24").Value = 11.76199650858811
$ws.Range("F24").Value = 26.74972159455259
$ws.Range("G24").Value = 23.04021129152765
$ws.Range("H24").Value = 12.80100387155337
$ws.Range("I24").Value = 17.54970647902338
$ws.Range("J24").Value = 11.06740321829899
$ws.Range("M24").Value = 17.81710583586217
$ws.Range("O24").Value = 18.72548171869391
$ws.Range("B25").Value = 13.87049492967503
$ws.Range("C25").Value = 9.561582368614641
$ws.Range("D25").Value = 11.70629246321568
$ws.Range("F25").Value = 26.87392366772989
$ws.Range("G25").Value = 23.17134862623188
$ws.Range("H25").Value = 12.89749637646021
$ws.Range("I25").Value = 17.79085503259386
$ws.Range("J25").Value = 11.12900022184246
$ws.Range("M25").Value = 17.34317179109349
$ws.Range("O25").Value = 18.88103600722202
